# "Generate Report for Handoff"
#
# The localization status report moves from "In Translation" to
# "Ready for handoff": the per-language status cells change, the
# "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
# advance a bit, and the Status column widens to fit the longer text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value     = "Ready for handoff"
$dede.Range("C2").Value     = "Ready for handoff"

# --- Timestamps advance a little as the handoff package is generated ---
$overview.Range("G2").Value = "2016-09-06 00:44:38"
$dede.Range("H2").Value     = "2016-09-06 00:44:38"
$zhcn.Range("H2").Value     = "2016-09-06 00:44:34"

# --- Widen the Status columns so the longer text still fits ---
# (same target character width applied on Overview!E:F and the Status
# column, C, on each per-language sheet)
$newStatusColWidth = 16.333333333333332
$overview.Columns.Item(5).ColumnWidth = $newStatusColWidth
$overview.Columns.Item(6).ColumnWidth = $newStatusColWidth
$zhcn.Columns.Item(3).ColumnWidth     = $newStatusColWidth
$dede.Columns.Item(3).ColumnWidth     = $newStatusColWidth
